# "new code 11 added for plotting"
# Fill column B (the "done / plotted" checkmark column) down from the last
# populated cell (B11) through the end of the data in column A (row 61),
# then move the active selection to D3 (no more frozen/scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1:B11 already contain the "✅" marker; extend it down through row 61
# so every row in the A column gets the same marker in column B.
$ws.Range("B11:B61").FillDown() | Out-Null

# FillDown correctly propagates the value/formula but not the cell style,
# so copy the source cell's formatting (centered alignment, style index 1)
# onto the newly filled cells without touching their values.
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B12:B61").PasteSpecial(-4122) | Out-Null

# Update the visible selection/active cell to D3.
$ws.Range("D3").Select() | Out-Null
